# Auto-generated edit script applying scheduled-runner market data updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H18").Value = 24949.75
$ws.Range("J18").Value = 48999.5
$ws.Range("L18").Value = 48999.5
$ws.Range("N18").Value = -49567.5

$ws.Range("H123").Value = 65166.668
$ws.Range("J123").Value = 65166.668
$ws.Range("L123").Value = 65166.668
$ws.Range("N123").Value = -74966.66800000001

$ws.Range("H124").Value = 80000
$ws.Range("J124").Value = 80000
$ws.Range("L124").Value = 80000
$ws.Range("N124").Value = -89820

$ws.Range("H132").Value = 5894.857
$ws.Range("I132").Value = 1844.0769
$ws.Range("K132").Value = 5532.2307
$ws.Range("M132").Value = -3002.2307

$ws.Range("H137").Value = 38548.258
$ws.Range("I137").Value = 1491.5217
$ws.Range("J137").Value = 251624.5
$ws.Range("K137").Value = 4474.5651
$ws.Range("L137").Value = 754873.5
$ws.Range("M137").Value = -1924.5651
$ws.Range("N137").Value = -759973.5

$ws.Range("H138").Value = 1748.3914
$ws.Range("I138").Value = 1323.8889
$ws.Range("J138").Value = 3276.6
$ws.Range("K138").Value = 3971.6667
$ws.Range("L138").Value = 9829.799999999999
$ws.Range("M138").Value = 1168.3333
$ws.Range("N138").Value = -20109.8

$ws = $wb.Worksheets("ARM")
$ws.Range("H2").Value = 2672.7727
$ws.Range("I2").Value = 2451.4443
$ws.Range("K2").Value = 2451.4443
$ws.Range("M2").Value = -2338.4443

$ws.Range("H45").Value = 23279.268
$ws.Range("I45").Value = 24983
$ws.Range("K45").Value = 24983
$ws.Range("M45").Value = -24606

$ws.Range("H116").Value = 2672.7727
$ws.Range("I116").Value = 2451.4443
$ws.Range("K116").Value = 2451.4443
$ws.Range("M116").Value = -157.4443000000001

$ws.Range("H133").Value = 90000
$ws.Range("J133").Value = 90000
$ws.Range("L133").Value = 90000
$ws.Range("N133").Value = -95060

$ws = $wb.Worksheets("BSM")
$ws.Range("H3").Value = 2672.7727
$ws.Range("I3").Value = 2451.4443
$ws.Range("K3").Value = 2451.4443
$ws.Range("M3").Value = -2337.4443

$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H86").Value = 1889.6072
$ws.Range("I86").Value = 1583
$ws.Range("J86").Value = 3013.8333
$ws.Range("K86").Value = 1583
$ws.Range("L86").Value = 3013.8333
$ws.Range("M86").Value = -460
$ws.Range("N86").Value = -5259.8333

$ws.Range("H87").Value = 30354
$ws.Range("J87").Value = 30354
$ws.Range("L87").Value = 30354
$ws.Range("N87").Value = -32850

$ws.Range("H89").Value = 1889.6072
$ws.Range("I89").Value = 1583
$ws.Range("J89").Value = 3013.8333
$ws.Range("K89").Value = 7915
$ws.Range("L89").Value = 15069.1665
$ws.Range("M89").Value = -2299
$ws.Range("N89").Value = -26301.1665

$ws.Range("H90").Value = 30354
$ws.Range("J90").Value = 30354
$ws.Range("L90").Value = 91062
$ws.Range("N90").Value = -103542

$ws.Range("H107").Value = 18203.818
$ws.Range("I107").Value = 6521.7144
$ws.Range("K107").Value = 6521.7144
$ws.Range("M107").Value = -4601.7144

$ws = $wb.Worksheets("CRP")
$ws.Range("H38").Value = 4085.2856
$ws.Range("I38").Value = 1619.8823
$ws.Range("J38").Value = 7895.4546
$ws.Range("K38").Value = 1619.8823
$ws.Range("L38").Value = 7895.4546
$ws.Range("M38").Value = -1242.8823
$ws.Range("N38").Value = -8649.454600000001

$ws.Range("H46").Value = 4085.2856
$ws.Range("I46").Value = 1619.8823
$ws.Range("J46").Value = 7895.4546
$ws.Range("K46").Value = 1619.8823
$ws.Range("L46").Value = 7895.4546
$ws.Range("M46").Value = -1408.8823
$ws.Range("N46").Value = -8317.454600000001

$ws.Range("H50").Value = 10499.2
$ws.Range("J50").Value = 10499.2
$ws.Range("L50").Value = 10499.2
$ws.Range("N50").Value = -11749.2

$ws.Range("H51").Value = 11666
$ws.Range("J51").Value = 11666
$ws.Range("L51").Value = 11666
$ws.Range("N51").Value = -13138

$ws.Range("H60").Value = 7909
$ws.Range("J60").Value = 7856.857
$ws.Range("L60").Value = 7856.857
$ws.Range("N60").Value = -8878.857

$ws.Range("H61").Value = 11666
$ws.Range("J61").Value = 11666
$ws.Range("L61").Value = 11666
$ws.Range("N61").Value = -12362

$ws.Range("H62").Value = 3659.8
$ws.Range("I62").Value = 3566.3333
$ws.Range("J62").Value = 3800
$ws.Range("K62").Value = 3566.3333
$ws.Range("L62").Value = 3800
$ws.Range("M62").Value = -2942.3333
$ws.Range("N62").Value = -5048

$ws.Range("H65").Value = 3659.8
$ws.Range("I65").Value = 3566.3333
$ws.Range("J65").Value = 3800
$ws.Range("K65").Value = 17831.6665
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = -14711.6665
$ws.Range("N65").Value = -25240

$ws.Range("H99").Value = 15898.25
$ws.Range("I99").Value = 4500
$ws.Range("J99").Value = 27296.5
$ws.Range("K99").Value = 4500
$ws.Range("L99").Value = 27296.5
$ws.Range("M99").Value = -3002
$ws.Range("N99").Value = -30292.5

$ws.Range("H126").Value = 15898.25
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 27296.5
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 81889.5
$ws.Range("M126").Value = -11030
$ws.Range("N126").Value = -86829.5

$ws.Range("H134").Value = 1999.875
$ws.Range("I134").Value = 1858.3704
$ws.Range("K134").Value = 5575.1112
$ws.Range("M134").Value = -3040.1112

$ws = $wb.Worksheets("CUL")
$ws.Range("H97").Value = 1423.1
$ws.Range("I97").Value = 1199.6666
$ws.Range("J97").Value = 1518.8572
$ws.Range("K97").Value = 3598.9998
$ws.Range("L97").Value = 4556.571599999999
$ws.Range("M97").Value = -3102.9998
$ws.Range("N97").Value = -5548.571599999999

$ws.Range("H129").Value = 98358.62
$ws.Range("J129").Value = 3462.8333
$ws.Range("L129").Value = 10388.4999
$ws.Range("N129").Value = -20388.4999

$ws.Range("H131").Value = 1700.6957
$ws.Range("I131").Value = 1641.2222
$ws.Range("J131").Value = 1738.9286
$ws.Range("K131").Value = 4923.6666
$ws.Range("L131").Value = 5216.7858
$ws.Range("M131").Value = 116.3334000000004
$ws.Range("N131").Value = -15296.7858

$ws = $wb.Worksheets("GSM")
$ws.Range("H113").Value = 50002584
$ws.Range("I113").Value = 55557764
$ws.Range("K113").Value = 55557764
$ws.Range("M113").Value = -55555594

$ws.Range("H126").Value = 3364.5715
$ws.Range("I126").Value = 2499
$ws.Range("K126").Value = 7497
$ws.Range("M126").Value = -5027

$ws.Range("H136").Value = 51662.5
$ws.Range("J136").Value = 51662.5
$ws.Range("L136").Value = 154987.5
$ws.Range("N136").Value = -160087.5

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 41372
$ws.Range("I7").Value = 70001
$ws.Range("K7").Value = 70001
$ws.Range("M7").Value = -69889

$ws.Range("H22").Value = 1678.2142
$ws.Range("I22").Value = 1677.2222
$ws.Range("J22").Value = 1680
$ws.Range("K22").Value = 1677.2222
$ws.Range("L22").Value = 1680
$ws.Range("M22").Value = -1382.2222
$ws.Range("N22").Value = -2270

$ws.Range("H27").Value = 1678.2142
$ws.Range("I27").Value = 1677.2222
$ws.Range("J27").Value = 1680
$ws.Range("K27").Value = 1677.2222
$ws.Range("L27").Value = 1680
$ws.Range("M27").Value = -1570.2222
$ws.Range("N27").Value = -1894

$ws.Range("H46").Value = 2615.1177
$ws.Range("I46").Value = 2176
$ws.Range("K46").Value = 2176
$ws.Range("M46").Value = -1988

$ws.Range("H122").Value = 3987.1428
$ws.Range("I122").Value = 3475
$ws.Range("K122").Value = 10425
$ws.Range("M122").Value = -7975

$ws.Range("H126").Value = 41372
$ws.Range("I126").Value = 70001
$ws.Range("K126").Value = 210003
$ws.Range("M126").Value = -207533

$ws.Range("H132").Value = 2709.5715
$ws.Range("I132").Value = 2513.84
$ws.Range("J132").Value = 3198.9
$ws.Range("K132").Value = 7541.52
$ws.Range("L132").Value = 9596.700000000001
$ws.Range("M132").Value = -5011.52
$ws.Range("N132").Value = -14656.7

$ws.Range("H133").Value = 79924.5
$ws.Range("J133").Value = 79924.5
$ws.Range("L133").Value = 79924.5
$ws.Range("N133").Value = -84984.5

$ws = $wb.Worksheets("WVR")
$ws.Range("H15").Value = 11982.5
$ws.Range("J15").Value = 11982.5
$ws.Range("L15").Value = 11982.5
$ws.Range("N15").Value = -12558.5

$ws.Range("H54").Value = 40022.668
$ws.Range("J54").Value = 49999
$ws.Range("L54").Value = 49999
$ws.Range("N54").Value = -51039

$ws.Range("H122").Value = 3848.4443
$ws.Range("I122").Value = 1846.6666
$ws.Range("J122").Value = 4849.3335
$ws.Range("K122").Value = 5539.9998
$ws.Range("L122").Value = 14548.0005
$ws.Range("M122").Value = -3089.9998
$ws.Range("N122").Value = -19448.0005

$ws.Range("H126").Value = 2680.4285
$ws.Range("J126").Value = 4500
$ws.Range("L126").Value = 13500
$ws.Range("N126").Value = -18440

$ws.Range("H132").Value = 6442.795
$ws.Range("I132").Value = 6621.857
$ws.Range("K132").Value = 19865.571
$ws.Range("M132").Value = -17335.571

$ws.Range("H136").Value = 3647.0667
$ws.Range("I136").Value = 2142.25
$ws.Range("J136").Value = 9666.333000000001
$ws.Range("K136").Value = 6426.75
$ws.Range("L136").Value = 28998.999
$ws.Range("M136").Value = -3876.75
$ws.Range("N136").Value = -34098.999
